# Add "UK" Test Data sheet, modeled on the existing "Poland" sheet
# (same layout, plus the two extra product rows P32AR / P32DR that
# Poland's market doesn't have), matching the commit "Added Test Data
# for UK Market".

$wb = $excel.ActiveWorkbook

# Duplicate the "Poland" sheet (closest template: same row layout,
# same row-5 height, same merged cells) and place the copy right after it.
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)
$ws = $wb.Worksheets.Item("Poland (2)")
$ws.Name = "UK"

# Poland's sheet is missing the "P32AR" / "P32DR" product rows that the
# fuller market sheets (e.g. Germany) include - insert them in place,
# right before "PR1DS".
$ws.Rows("16:17").Insert()

# Copy formatting down from the row below onto the two new rows, then
# set their text.
$ws.Range("A18").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$ws.Range("A16").Value = "P32AR"
$ws.Range("A17").Value = "P32DR"

# Market-specific header values (ticket reference first, then the
# market name, so new shared-string entries land in that order).
$ws.Range("B4").Value = "NGC-2741/T3365"
$ws.Range("B2").Value = "UK Market"

# Column B is narrower on the UK sheet than on Poland's.
$ws.Columns("B").ColumnWidth = 15

# Leave the cursor on the input-value cell, like the source edit.
[void]$ws.Range("B4").Select()
